# Update the two-digit division worksheet numbers to the new set of
# problems. The worksheet is a single 5-column table; each "problem" row
# (1, 5, 9, 13, 17) is followed by three blank rows used for working out
# the answer. We replace the text of each of the 25 populated cells with
# its new value, scoped via Find.Execute on each cell's Range so that
# existing run formatting (font/size) is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row, Col, OldText, NewText
$replacements = @(
    @(1, 1, "99÷9=", "77÷5="),
    @(1, 2, "45÷5=", "84÷2="),
    @(1, 3, "94÷7=", "67÷5="),
    @(1, 4, "64÷3=", "74÷5="),
    @(1, 5, "32÷5=", "94÷2="),

    @(5, 1, "15÷4=", "94÷8="),
    @(5, 2, "24÷8=", "14÷9="),
    @(5, 3, "74÷7=", "19÷4="),
    @(5, 4, "10÷6=", "17÷9="),
    @(5, 5, "41÷7=", "23÷7="),

    @(9, 1, "18÷7=", "49÷9="),
    @(9, 2, "37÷2=", "10÷7="),
    @(9, 3, "64÷3=", "42÷2="),
    @(9, 4, "86÷3=", "17÷2="),
    @(9, 5, "30÷8=", "92÷6="),

    @(13, 1, "72÷2=", "93÷5="),
    @(13, 2, "81÷9=", "92÷3="),
    @(13, 3, "67÷4=", "42÷2="),
    @(13, 4, "30÷3=", "17÷2="),
    @(13, 5, "90÷7=", "59÷6="),

    @(17, 1, "56÷4=", "14÷6="),
    @(17, 2, "84÷4=", "42÷2="),
    @(17, 3, "21÷7=", "53÷5="),
    @(17, 4, "89÷8=", "26÷8="),
    @(17, 5, "57÷7=", "22÷9=")
)

foreach ($item in $replacements) {
    $r = $item[0]
    $c = $item[1]
    $old = $item[2]
    $new = $item[3]
    $cell = $t.Cell($r, $c)
    $rng = $cell.Range
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 0, $false, $new, 1)
}
